$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values - Invoice No, Task Type, Tariff Fee, Gross Fee, VAT %, VAT Amount, Net Fee,
# Case Details, Submission Date, Invoice Date, Payment Status
$ws.Range("A2").Value = "avv1"
$ws.Range("B2").Value = "dee"
$ws.Range("C2").Value = 1000
$ws.Range("D2").Value = 100
$ws.Range("E2").Value = 20
$ws.Range("F2").Value = 20
$ws.Range("G2").Value = 80
$ws.Range("H2").Value = "denememe"

$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "11.02.2025"

$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "12.02.2025"

$ws.Range("K2").Value = "Paid"
